$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 278065.38
$ws.Range("J17").Value = 285938.66
$ws.Range("L17").Value = 857815.98
$ws.Range("N17").Value = -858151.98
$ws.Range("H86").Value = 4131.9
$ws.Range("J86").Value = 4338.8
$ws.Range("L86").Value = 4338.8
$ws.Range("N86").Value = -6584.8
$ws.Range("H89").Value = 4131.9
$ws.Range("J89").Value = 4338.8
$ws.Range("L89").Value = 21694
$ws.Range("N89").Value = -32926
$ws.Range("H98").Value = 2970074
$ws.Range("I98").Value = 2934274.8
$ws.Range("K98").Value = 2934274.8
$ws.Range("M98").Value = -2932776.8
$ws.Range("H111").Value = 1199.5
$ws.Range("I111").Value = 1199.5
$ws.Range("K111").Value = 3598.5
$ws.Range("M111").Value = -531.5
$ws.Range("H122").Value = 2970074
$ws.Range("I122").Value = 2934274.8
$ws.Range("K122").Value = 8802824.399999999
$ws.Range("M122").Value = -8800374.399999999
$ws.Range("H131").Value = 21812.166
$ws.Range("I131").Value = 22574.7
$ws.Range("K131").Value = 67724.10000000001
$ws.Range("M131").Value = -62684.10000000001
$ws.Range("H138").Value = 2146.4387
$ws.Range("I138").Value = 968.1579
$ws.Range("J138").Value = 2892.6833
$ws.Range("K138").Value = 2904.4737
$ws.Range("L138").Value = 8678.0499
$ws.Range("M138").Value = 2235.5263
$ws.Range("N138").Value = -18958.0499

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 34588
$ws.Range("J76").Value = 34588
$ws.Range("L76").Value = 34588
$ws.Range("N76").Value = -35218
$ws.Range("H79").Value = 34588
$ws.Range("J79").Value = 34588
$ws.Range("L79").Value = 34588
$ws.Range("N79").Value = -36772

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 279432.62
$ws.Range("I31").Value = 371896.56
$ws.Range("J31").Value = 2040.8889
$ws.Range("K31").Value = 371896.56
$ws.Range("L31").Value = 2040.8889
$ws.Range("M31").Value = -371601.56
$ws.Range("N31").Value = -2630.8889
$ws.Range("H34").Value = 279432.62
$ws.Range("I34").Value = 371896.56
$ws.Range("J34").Value = 2040.8889
$ws.Range("K34").Value = 371896.56
$ws.Range("L34").Value = 2040.8889
$ws.Range("M34").Value = -371694.56
$ws.Range("N34").Value = -2444.8889
$ws.Range("H58").Value = 2609.25
$ws.Range("I58").Value = 2357.2942
$ws.Range("J58").Value = 4037
$ws.Range("K58").Value = 2357.2942
$ws.Range("L58").Value = 4037
$ws.Range("M58").Value = -2154.2942
$ws.Range("N58").Value = -4443
$ws.Range("H99").Value = 487880.2
$ws.Range("I99").Value = 838790.9399999999
$ws.Range("J99").Value = 19999.223
$ws.Range("K99").Value = 838790.9399999999
$ws.Range("L99").Value = 19999.223
$ws.Range("M99").Value = -837292.9399999999
$ws.Range("N99").Value = -22995.223
$ws.Range("H122").Value = 6446.6
$ws.Range("I122").Value = 6447.5
$ws.Range("J122").Value = 6443
$ws.Range("K122").Value = 19342.5
$ws.Range("L122").Value = 19329
$ws.Range("M122").Value = -16892.5
$ws.Range("N122").Value = -24229
$ws.Range("H126").Value = 487880.2
$ws.Range("I126").Value = 838790.9399999999
$ws.Range("J126").Value = 19999.223
$ws.Range("K126").Value = 2516372.82
$ws.Range("L126").Value = 59997.66900000001
$ws.Range("M126").Value = -2513902.82
$ws.Range("N126").Value = -64937.66900000001
$ws.Range("H132").Value = 3267.3572
$ws.Range("I132").Value = 3250.2354
$ws.Range("K132").Value = 9750.706200000001
$ws.Range("M132").Value = -7220.706200000001
$ws.Range("H134").Value = 5368.2285
$ws.Range("I134").Value = 5838.357
$ws.Range("K134").Value = 17515.071
$ws.Range("M134").Value = -14980.071
$ws.Range("H136").Value = 2609.25
$ws.Range("I136").Value = 2357.2942
$ws.Range("J136").Value = 4037
$ws.Range("K136").Value = 7071.882599999999
$ws.Range("L136").Value = 12111
$ws.Range("M136").Value = -4521.882599999999
$ws.Range("N136").Value = -17211

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 488.125
$ws.Range("I60").Value = 434.5
$ws.Range("J60").Value = 649
$ws.Range("K60").Value = 1303.5
$ws.Range("L60").Value = 1947
$ws.Range("M60").Value = -1052.5
$ws.Range("N60").Value = -2449

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9010.556
$ws.Range("J113").Value = 4413.7144
$ws.Range("L113").Value = 4413.7144
$ws.Range("N113").Value = -8753.714400000001
$ws.Range("H122").Value = 8895.290000000001
$ws.Range("I122").Value = 7778.25
$ws.Range("K122").Value = 23334.75
$ws.Range("M122").Value = -20884.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6932.517
$ws.Range("J7").Value = 7275
$ws.Range("L7").Value = 7275
$ws.Range("N7").Value = -7499
$ws.Range("H40").Value = 13932.523
$ws.Range("I40").Value = 14788.053
$ws.Range("K40").Value = 14788.053
$ws.Range("M40").Value = -14652.053
$ws.Range("H46").Value = 2198.3
$ws.Range("I46").Value = 2198
$ws.Range("J46").Value = 2198.6
$ws.Range("K46").Value = 2198
$ws.Range("L46").Value = 2198.6
$ws.Range("M46").Value = -2010
$ws.Range("N46").Value = -2574.6
$ws.Range("H122").Value = 3049.2778
$ws.Range("I122").Value = 2780.7273
$ws.Range("K122").Value = 8342.1819
$ws.Range("M122").Value = -5892.1819
$ws.Range("H126").Value = 6932.517
$ws.Range("J126").Value = 7275
$ws.Range("L126").Value = 21825
$ws.Range("N126").Value = -26765
$ws.Range("H132").Value = 7486.476
$ws.Range("I132").Value = 8891.200000000001
$ws.Range("J132").Value = 3974.6667
$ws.Range("K132").Value = 26673.6
$ws.Range("L132").Value = 11924.0001
$ws.Range("M132").Value = -24143.6
$ws.Range("N132").Value = -16984.0001
$ws.Range("H136").Value = 1372.9333
$ws.Range("I136").Value = 1249.5834
$ws.Range("J136").Value = 1866.3334
$ws.Range("K136").Value = 3748.7502
$ws.Range("L136").Value = 5599.0002
$ws.Range("M136").Value = -1198.7502
$ws.Range("N136").Value = -10699.0002

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 22857
$ws.Range("J14").Value = 18666.334
$ws.Range("L14").Value = 18666.334
$ws.Range("N14").Value = -19002.334
$ws.Range("H122").Value = 11865.6455
$ws.Range("I122").Value = 2965.087
$ws.Range("K122").Value = 8895.261
$ws.Range("M122").Value = -6445.261
$ws.Range("H126").Value = 3183.3076
$ws.Range("J126").Value = 2800.8
$ws.Range("L126").Value = 8402.400000000001
$ws.Range("N126").Value = -13342.4
$ws.Range("H136").Value = 657642.2
$ws.Range("I136").Value = 910552.8
$ws.Range("J136").Value = 101238.8
$ws.Range("K136").Value = 2731658.4
$ws.Range("L136").Value = 303716.4
$ws.Range("M136").Value = -2729108.4
$ws.Range("N136").Value = -308816.4
